# Update to framework 7.1
#
# Adds four new "Model_*" formula-library rows to the Library_Formula sheet,
# matching the look/feel of the existing rows (same Action/Library columns,
# String in/out types), and updates the saved selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

$newNames = @("Model_1_Corporate", "Model_2_Corporate", "Model_1_Retail", "Model_2_Retail")
$lastRow = 106

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $row = $lastRow + 1 + $i

    # Duplicate the last existing data row (same styling on every column,
    # including the blank D/G columns) onto the new row.
    $ws.Rows("$lastRow`:$lastRow").Copy()
    $ws.Rows("$row`:$row").Insert(-4121)

    # Overwrite the Formula-Name cell with the new value.
    $ws.Cells.Item($row, 3).Value = $newNames[$i]

    # The Formula-Name column normally carries the "default" cell style
    # (same as column A/E/F) rather than the slightly different style used
    # on row 106's Formula-Name cell, so re-apply that default format.
    $ws.Cells.Item(2, 3).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Match the sheet's saved selection/state.
$ws.Activate()
$ws.Range("D111").Select()
